$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the provider slug and the "benefits" key to "rates" (william russell rates update)
$ws.Range("A2").Value = "william_russell"
$ws.Range("J2").Value = "rates"

# Move the active selection as left by the author
[void]$ws.Range("A6").Select()
